$d = $word.ActiveDocument

$hit = $d.Content
$hit.Find.Execute("aaa", $false, $false, $false, $false, $false,
                   $true, 1, $false, "", 0) | Out-Null
$hit.Collapse(1)  # wdCollapseStart = 1

$hit.InsertBefore("c")
$hit.Collapse(1)
$hit.InsertBefore("ff")
$hit.Collapse(1)
$hit.InsertBefore("d")
